# edit.ps1
# Applies the three-part change described by the diff:
#  1. Slide 16: title run "Τεμαχισμός Συμβολοσειράς" is split into two
#     runs: "Διαμέριση" + " Συμβολοσειράς" (the word "Τεμαχισμός" is
#     replaced by "Διαμέριση").
#  2. Slide 19: body placeholder (2nd shape) is resized/repositioned.
#  3. Slide 19: body placeholder text about the `in` keyword is reworded.

$p = $ppt.ActivePresentation

# EMU per point, used throughout for Left/Top/Width/Height (which the
# PowerPoint object model exposes in points even though OOXML stores EMU).
# A tiny +0.5 EMU nudge before dividing compensates for truncation that
# happens when the point value is converted back to EMU internally.
$EMU_PER_PT = 12700

function EmuToPt($emu) {
    return ($emu + 0.5) / $EMU_PER_PT
}

# ---------------------------------------------------------------------
# 1) Slide 16 - title placeholder: "Τεμαχισμός" -> "Διαμέριση"
# ---------------------------------------------------------------------
$slide16 = $p.Slides.Item(16)
$titleShape = $slide16.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange

# The run currently reads "Τεμαχισμός Συμβολοσειράς" (24 chars); the
# first 10 characters ("Τεμαχισμός") become "Διαμέριση" while the rest
# of the run (" Συμβολοσειράς") keeps its own separate run. Setting the
# .Text of a Characters() sub-range naturally splits the parent run into
# multiple runs, each preserving the original run formatting.
$titleFirstWord = $titleRange.Characters(1, 10)
$titleFirstWord.Text = "Διαμέριση"

# ---------------------------------------------------------------------
# 2) Slide 19 - body placeholder: reposition/resize + reword
# ---------------------------------------------------------------------
$slide19 = $p.Slides.Item(19)
$bodyShape = $slide19.Shapes.Item(2)

$bodyShape.Left = EmuToPt 1155699
$bodyShape.Top = EmuToPt 2921548
$bodyShape.Width = EmuToPt 7246179
$bodyShape.Height = EmuToPt 5702399

$bodyRange = $bodyShape.TextFrame.TextRange

# First paragraph runs: "Η λέξη-κλειδί " (14) + "in" (2) + the long
# explanation run (103 chars) that needs to be reworded.
$explanationRun = $bodyRange.Characters(17, 103)
$explanationRun.Text = " μπορεί επίσης να χρησιμοποιηθεί για να ελέγξει εάν μια συμβολοσειρά αποτελεί τμήμα μιας άλλης συμβολοσειράς"

Write-Host "edit complete"
